# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-03 12:32:22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded-by email list reorders / updates ---
$ws.Range("G2").Value  = "rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G18").Value = "shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G40").Value = "shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G58").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G62").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G80").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G84").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G96").Value = "Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G98").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G106").Value = "nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G118").Value = "Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G120").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G128").Value = "nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G134").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G150").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G172").Value = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Student headcount updates tied to the Recorded-by changes above ---
$ws.Range("H150").Value = "69/224"
$ws.Range("H172").Value = "34/226"

# --- Class Statistics block (K/L columns) ---
$ws.Range("L6").Value  = 28
$ws.Range("L7").Value  = 10
$ws.Range("L9").Value  = "15.9%"
$ws.Range("L10").Value = "29.6%"

# --- Subgroup statistics block (K:S columns) ---
# Row 16 (Year 2 / A2)
$ws.Range("O16").Value = 3
$ws.Range("P16").Value = 2
$ws.Range("R16").Value = "13.6%"
$ws.Range("S16").Value = "30.0%"

# Row 17 (Year 2 / A3)
$ws.Range("O17").Value = 4
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "18.2%"
$ws.Range("S17").Value = "37.4%"

# Row 21 (Year 2 / B3)
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 0
$ws.Range("R21").Value = "18.2%"
$ws.Range("S21").Value = "24.1%"

# Row 22 (Year 2 / B4)
$ws.Range("S22").Value = "15.0%"

# --- Rows that flipped from "Not Recorded" to "Recorded" (also re-colors A:I like other Recorded rows) ---
$ws.Range("A2:I2").Copy()
$ws.Range("A36:I36").PasteSpecial(-4122)
$ws.Range("G36").Value = "Amr-Saeed@med.asu.edu.eg"
$ws.Range("H36").Value = "53/217"
$ws.Range("I36").Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A54:I54").PasteSpecial(-4122)
$ws.Range("G54").Value = "Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("H54").Value = "41/220"
$ws.Range("I54").Value = "Recorded"

$ws.Range("A2:I2").Copy()
$ws.Range("A146:I146").PasteSpecial(-4122)
$ws.Range("G146").Value = "Amr-Saeed@med.asu.edu.eg"
$ws.Range("H146").Value = "57/224"
$ws.Range("I146").Value = "Recorded"
